# #4801_PreRegister - payment_export/export_template/xlsx_report_register_icash_smart.xlsx
# Adds a new "Run Date" header row (row 3) above the existing column-header
# row, styled like the other yellow report headers, plus a light-yellow
# highlighted input cell next to it, and selects C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 / A3: "Run Date" label -----------------------------------------
# Reuse the existing bold "TH SarabunPSK" header font + thin box border
# (fontId 2 / borderId 2, already used by B5) by copy/pasting formats from
# B5, then re-color the fill to the report's existing yellow accent and
# left-align the text instead of the centered original.
$ws.Range("B5").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = "Run Date"
$ws.Range("A3").Interior.Color = 65535        # RGB(255,255,0) -> FFFFFF00 (existing yellow fill)
$ws.Range("A3").HorizontalAlignment = -4131   # xlLeft
$ws.Range("A3").VerticalAlignment = -4107     # xlBottom (default -> keeps alignment minimal)

# --- Row 3 / B3: light-yellow highlighted (empty) input cell --------------
$ws.Range("B3").Interior.Color = 10092543     # RGB(255,255,153) -> FFFFFF99

# --- Row height for the new row --------------------------------------------
$ws.Rows.Item(3).RowHeight = 18

# --- Selection matches the authored sheet view ------------------------------
$ws.Range("C3").Select()
